$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the header row (row 1) and delete it entirely, shifting data up.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
